# Refresh the cryptos list: update Price (D) and Volume(1h) (E) columns
# for rows 2-51. Values are written as text (leading apostrophe) to match
# the source data's text representation (prices such as "29.175.96" are
# not valid numbers, and the sheet stores them - and the percentages - as
# plain text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{ Row = 2;  D = "29.175.96";    E = "  -2.05%  " }
  @{ Row = 3;  D = "1.855.15";     E = "  -0.95%  " }
  @{ Row = 4;  D = "0.9992";       E = "  -0.29%  " }
  @{ Row = 5;  D = "238.21";       E = "  -1.50%  " }
  @{ Row = 6;  D = "0.6914";       E = "  -3.76%  " }
  @{ Row = 7;  D = "0.9996";       E = "  -0.18%  " }
  @{ Row = 8;  D = "0.07720";      E = "  +2.46%  " }
  @{ Row = 9;  D = "0.3052";       E = "  -2.96%  " }
  @{ Row = 10; D = "23.27";        E = "  -5.33%  " }
  @{ Row = 11; D = "0.08018";      E = "  -2.06%  " }
  @{ Row = 12; D = "1.856.51";     E = "  -1.55%  " }
  @{ Row = 13;                     E = "  -2.78%  " }
  @{ Row = 14; D = "5.204";        E = "  -2.46%  " }
  @{ Row = 15; D = "89.39";        E = "  -3.28%  " }
  @{ Row = 16; D = "29.206.60";    E = "  -1.80%  " }
  @{ Row = 17; D = "5.744";        E = "  -4.56%  " }
  @{ Row = 18; D = "0.000007805";  E = "  -1.53%  " }
  @{ Row = 19; D = "13.24";        E = "  -1.69%  " }
  @{ Row = 20; D = "235.20";       E = "  -4.52%  " }
  @{ Row = 21; D = "0.9999";       E = "  +0.00%  " }
  @{ Row = 22; D = "2.105.29";     E = "  -0.38%  " }
  @{ Row = 23; D = "0.9991";       E = "  -0.30%  " }
  @{ Row = 24; D = "7.471";        E = "  -3.15%  " }
  @{ Row = 25; D = "161.98";       E = "  -1.09%  " }
  @{ Row = 26; D = "8.969";        E = "  -2.34%  " }
  @{ Row = 27; D = "0.1445";       E = "  -3.60%  " }
  @{ Row = 28;                     E = "  -2.64%  " }
  @{ Row = 29; D = "1.961";        E = "  -2.09%  " }
  @{ Row = 30; D = "1.405";        E = "  -1.62%  " }
  @{ Row = 31; D = "4.519";        E = "  -0.70%  " }
  @{ Row = 32;                     E = "  -2.33%  " }
  @{ Row = 33; D = "4.019";        E = "  -3.79%  " }
  @{ Row = 34; D = "0.05182";      E = "  -4.83%  " }
  @{ Row = 35; D = "1.186";        E = "  -3.05%  " }
  @{ Row = 36; D = "0.7042";       E = "  -4.50%  " }
  @{ Row = 37; D = "1.006";        E = "  +0.72%  " }
  @{ Row = 38;                     E = "  -1.20%  " }
  @{ Row = 39;                     E = "  -3.15%  " }
  @{ Row = 40; D = "2.679";        E = "  -1.97%  " }
  @{ Row = 41; D = "0.9266";       E = "  +4.02%  " }
  @{ Row = 42; D = "1.102.46";     E = "  +6.24%  " }
  @{ Row = 43; D = "5.961";        E = "  -0.55%  " }
  @{ Row = 44; D = "0.4287";       E = "  -3.74%  " }
  @{ Row = 45; D = "70.41";        E = "  -1.66%  " }
  @{ Row = 46; D = "0.9994";       E = "  -0.22%  " }
  @{ Row = 47; D = "102.47";       E = "  -1.04%  " }
  @{ Row = 48; D = "1.791";        E = "  -0.43%  " }
  @{ Row = 49; D = "2.002.55";     E = "  -0.44%  " }
  @{ Row = 50; D = "9.183";        E = "  -4.17%  " }
  @{ Row = 51; D = "7.007";        E = "  -6.13%  " }
)

foreach ($u in $updates) {
  $r = $u.Row
  if ($u.ContainsKey("D")) {
    $ws.Cells.Item($r, 4).Value = "'" + $u.D
  }
  if ($u.ContainsKey("E")) {
    $ws.Cells.Item($r, 5).Value = "'" + $u.E
  }
}
